# Adds a new "VATRIM" worksheet at the end of the workbook with a two-column
# time series (date labels in column A, VA_trim figures in column B), mirroring
# the existing PREVISION-style layout (bold orange header row).

$wb = $excel.ActiveWorkbook

$prev = $wb.Worksheets.Item("PREVISION")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)

# 1. Insert the new sheet right after the last existing sheet.
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "VATRIM"

# 2. Header row - literal text, styled like the other sheets' header rows.
$ws.Range("A1").Formula = '="time"'
$ws.Range("B1").Formula = '="VA_trim"'
$ws.Range("A1:B1").Copy()
$ws.Range("A1:B1").PasteSpecial(-4163)   # xlPasteValues - flatten formula to a literal value
$prev.Range("A1:B1").Copy()
$ws.Range("A1:B1").PasteSpecial(-4122)   # xlPasteFormats - reuse PREVISION's header look

# 3. Data rows: quarterly dates (text) in column A, VA_trim figures (numbers) in column B.
$dates = @("1997-01-01","1997-04-01","1997-07-01","1997-10-01","1998-01-01","1998-04-01","1998-07-01","1998-10-01","1999-01-01","1999-04-01","1999-07-01","1999-10-01","2000-01-01","2000-04-01","2000-07-01","2000-10-01","2001-01-01","2001-04-01","2001-07-01","2001-10-01","2002-01-01","2002-04-01","2002-07-01","2002-10-01","2003-01-01","2003-04-01","2003-07-01","2003-10-01","2004-01-01","2004-04-01","2004-07-01","2004-10-01","2005-01-01","2005-04-01","2005-07-01","2005-10-01","2006-01-01","2006-04-01","2006-07-01","2006-10-01","2007-01-01","2007-04-01","2007-07-01","2007-10-01","2008-01-01","2008-04-01","2008-07-01","2008-10-01","2009-01-01","2009-04-01","2009-07-01","2009-10-01","2010-01-01","2010-04-01","2010-07-01","2010-10-01","2011-01-01","2011-04-01","2011-07-01","2011-10-01","2012-01-01","2012-04-01","2012-07-01","2012-10-01","2013-01-01","2013-04-01","2013-07-01","2013-10-01","2014-01-01","2014-04-01","2014-07-01","2014-10-01","2015-01-01","2015-04-01","2015-07-01","2015-10-01","2016-01-01","2016-04-01","2016-07-01","2016-10-01","2017-01-01","2017-04-01","2017-07-01","2017-10-01","2018-01-01","2018-04-01","2018-07-01","2018-10-01","2019-01-01","2019-04-01","2019-07-01","2019-10-01","2020-01-01","2020-04-01","2020-07-01","2020-10-01","2021-01-01","2021-04-01","2021-07-01","2021-10-01")
$vals = @(132947.420754456,212179.739504027,212276.33848145,227150.506767315,141393.103621597,226795.429710672,226141.013334937,239206.415486733,151583.382586458,241196.029205423,240580.767817142,257163.48842593,162331.055571178,255636.177730075,254761.346080887,268794.802071233,172152.944149418,262959.217922613,259489.313219138,272553.853557768,165954.121361025,254451.023334684,251850.110041742,265605.637329344,165394.611610315,260338.040947399,259586.928674046,276279.52476645,169633.912671073,266450.157233896,266503.51816605,283239.162626905,181107.650782012,276683.096726239,274530.816452221,288409.591236635,181453.27855698,270137.047097937,263887.320242544,272851.771650393,169545.058154116,270276.397914076,268450.180533682,288675.05677442,171154.950219043,288532.41522993,289609.170758743,308698.106163459,183405.31431924,296357.683513839,291988.326762593,307280.861394363,168979.639599388,274326.632053725,270052.111531832,285782.401860167,153859.834863015,262279.881650588,262568.888808087,283600.019780695,158488.861968795,264177.013235019,257066.434451383,263471.781077842,147564.46110679,270161.370594591,270234.255710547,302423.029869192,174444.054917061,338611.306219553,342392.336499563,365941.859530217,203613.060426295,348369.318706765,341264.200852872,362592.580291987,197030.448484361,347853.535772949,345871.87410849,371745.201973878,211599.878486715,374923.023234774,374356.279947385,402250.239172714,218161.644542643,374051.031218892,370119.002327364,392365.714838474,224956.232518354,383870.997190286,379616.057817982,405535.752239579,225694.858612034,391619.225768974,389280.648419773,415879.343420858,232721.751457131,403904.474216792,401534.97246494,429647.176254566)

$rowCount = $dates.Length
for ($i = 0; $i -lt $rowCount; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 1).Formula = '="' + $dates[$i] + '"'
    $ws.Cells.Item($r, 2).Value = $vals[$i]
}

# Flatten the date formulas in column A down to literal text values (shared strings),
# matching how the rest of the workbook stores its date labels.
$dataRange = $ws.Range("A2:A" + ($rowCount + 1))
$dataRange.Copy()
$dataRange.PasteSpecial(-4163)           # xlPasteValues

$excel.CutCopyMode = 0
